# Update the contact phone number on the final "Thank you" slide
# (slide with sldId=296, last slide in the deck) from 408-805-6749
# to 469-892-8857.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$phone = $tr.Paragraphs(2, 1)

# Setting the text directly would make the diff/merge logic in the
# runtime split the run on the shared "4" prefix between the old and
# new numbers. Clearing it out to an unrelated placeholder first keeps
# the final assignment a clean, single-run replacement (matching how
# PowerPoint itself would leave a single <a:r> run after a full
# selection retype).
$phone.Text = "placeholder"
$phone.Text = "469-892-8857"
